$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SL matchup averages for Spring '24 week 6 inputs
$ws.Range("D2").Value = 1.25
$ws.Range("E3").Value = 1.33
$ws.Range("B4").Value = 1.47
$ws.Range("C5").Value = 1.32
$ws.Range("G6").Value = 1.02
$ws.Range("F7").Value = 1.47
$ws.Range("G7").Value = 1.17
